$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.72"
$ws.Range("E2").Value = "'2.30%"
$ws.Range("D3").Value = "'32.35"
$ws.Range("E3").Value = "'3.84%"
$ws.Range("D4").Value = "'4.971"
$ws.Range("E4").Value = "'0.77%"
$ws.Range("D5").Value = "'0.07723"
$ws.Range("E5").Value = "'5.39%"
$ws.Range("D6").Value = "'2.334"
$ws.Range("E6").Value = "'2.90%"
$ws.Range("D7").Value = "'7.944"
$ws.Range("E7").Value = "'2.75%"
$ws.Range("B8").Value = "'GateToken"
$ws.Range("C8").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.825"
$ws.Range("E8").Value = "'2.03%"
$ws.Range("B9").Value = "'MXToken"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9232"
$ws.Range("E9").Value = "'1.66%"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09937"
$ws.Range("E10").Value = "'22.96%"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1759"
$ws.Range("E11").Value = "'3.82%"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08431"
$ws.Range("E12").Value = "'3.56%"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03273"
$ws.Range("E13").Value = "'5.42%"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09871"
$ws.Range("E14").Value = "'-2.15%"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001474"
$ws.Range("E15").Value = "'-2.86%"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005713"
$ws.Range("E16").Value = "'-0.17%"
$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.513"
$ws.Range("E17").Value = "'0.84%"
$ws.Range("D18").Value = "'2.184"
$ws.Range("E18").Value = "'5.29%"
$ws.Range("D19").Value = "'0.3364"
$ws.Range("E19").Value = "'1.11%"
$ws.Range("E20").Value = "'2.53%"
$ws.Range("D21").Value = "'4.368"
$ws.Range("E21").Value = "'9.85%"
$ws.Range("D22").Value = "'0.2089"
$ws.Range("E22").Value = "'-0.77%"
$ws.Range("D23").Value = "'0.04537"
$ws.Range("E23").Value = "'-0.15%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'0.31%"
$ws.Range("D25").Value = "'0.004366"
$ws.Range("E25").Value = "'0.53%"
$ws.Range("D26").Value = "'0.0001293"
$ws.Range("E26").Value = "'-0.78%"
$ws.Range("D27").Value = "'0.0003376"
$ws.Range("E27").Value = "'-0.76%"
$ws.Range("D39").Value = "'0.01704"
$ws.Range("E39").Value = "'6.72%"
$ws.Range("D40").Value = "'0.04711"
$ws.Range("E40").Value = "'5.99%"
$ws.Range("D41").Value = "'0.007718"
$ws.Range("E41").Value = "'5.73%"
$ws.Range("D42").Value = "'0.009778"
$ws.Range("E42").Value = "'12.87%"
$ws.Range("D43").Value = "'0.1393"
$ws.Range("E43").Value = "'4.84%"
$ws.Range("D44").Value = "'0.002111"
$ws.Range("E44").Value = "'8.93%"
$ws.Range("D45").Value = "'0.009651"
$ws.Range("E45").Value = "'1.26%"
$ws.Range("D46").Value = "'0.00006082"
$ws.Range("E46").Value = "'1.98%"
$ws.Range("D47").Value = "'0.00000000746"
$ws.Range("E47").Value = "'-0.76%"
$ws.Range("D48").Value = "'2.551"
$ws.Range("E48").Value = "'13.83%"
$ws.Range("D49").Value = "'0.001989"
$ws.Range("E49").Value = "'-31.38%"
$ws.Range("D50").Value = "'0.00002089"
$ws.Range("E50").Value = "'-0.76%"
$ws.Range("D51").Value = "'0.0001989"
$ws.Range("E51").Value = "'-0.76%"
